# Daily attendance processing - 2025-10-14 08:25:40
# Rotates the "Recorded By" (column G) list of names/emails for each data row:
# the first entry in the comma-separated list is moved to the end of the list.
# Rows whose "Recorded By" value contains only a single entry are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ",\s*"
        if ($parts.Count -gt 1) {
            $rotated = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ", "
            $cell.Value = $rotated
        }
    }
}
